# Weekly update: a new week of price records is prepended to the data
# block for this market/product. Two new rows are inserted right after
# the existing row 373, pushing the rest of the table (previously rows
# 374:395) down to 376:397, and the two freshly inserted rows (374:375)
# are populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 374; everything below shifts down.
$ws.Rows("374:375").Insert()

# Row 374: Crespo record / Primera
$ws.Cells.Item(374, 1).Value = 11
$ws.Cells.Item(374, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(374, 3).Value = "Bíobío"
$ws.Cells.Item(374, 4).Value = 44826
$ws.Cells.Item(374, 5).Value = 8
$ws.Cells.Item(374, 6).Value = 100112006
$ws.Cells.Item(374, 7).Value = "Repollo"
$ws.Cells.Item(374, 8).Value = "Crespo record"
$ws.Cells.Item(374, 9).Value = "Primera"
$ws.Cells.Item(374, 10).Value = 1000
$ws.Cells.Item(374, 11).Value = 1500
$ws.Cells.Item(374, 12).Value = 1600
$ws.Cells.Item(374, 13).Value = 1550
$ws.Cells.Item(374, 14).Value = "$/unidad"
$ws.Cells.Item(374, 15).Value = "Región Metropolitana"
$ws.Cells.Item(374, 16).Value = 1550
$ws.Cells.Item(374, 17).Value = 1
$ws.Cells.Item(374, 18).Value = "Hortaliza"

# Row 375: Crespo record / Segunda
$ws.Cells.Item(375, 1).Value = 11
$ws.Cells.Item(375, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(375, 3).Value = "Bíobío"
$ws.Cells.Item(375, 4).Value = 44826
$ws.Cells.Item(375, 5).Value = 8
$ws.Cells.Item(375, 6).Value = 100112006
$ws.Cells.Item(375, 7).Value = "Repollo"
$ws.Cells.Item(375, 8).Value = "Crespo record"
$ws.Cells.Item(375, 9).Value = "Segunda"
$ws.Cells.Item(375, 10).Value = 500
$ws.Cells.Item(375, 11).Value = 1300
$ws.Cells.Item(375, 12).Value = 1300
$ws.Cells.Item(375, 13).Value = 1300
$ws.Cells.Item(375, 14).Value = "$/unidad"
$ws.Cells.Item(375, 15).Value = "Región Metropolitana"
$ws.Cells.Item(375, 16).Value = 1300
$ws.Cells.Item(375, 17).Value = 1
$ws.Cells.Item(375, 18).Value = "Hortaliza"
